{"js": "// Replace the multiplication-problem text in each table cell with the\n// updated operands/equation, keeping the surrounding run formatting\n// (font, size) untouched because we replace only the matched text range.\nconst replacements = [\n  [\"933\u00d73=\", \"361\u00d75=\"],\n  [\"234\u00d79=\", \"358\u00d79=\"],\n  [\"137\u00d73=\", \"382\u00d76=\"],\n  [\"335\u00d77=\", \"732\u00d74=\"],\n  [\"479\u00d75=\", \"577\u00d78=\"],\n  [\"550\u00d79=\", \"685\u00d73=\"],\n  [\"375\u00d76=\", \"219\u00d77=\"],\n  [\"680\u00d75=\", \"558\u00d79=\"],\n  [\"275\u00d74=\", \"195\u00d78=\"],\n  [\"694\u00d77=\", \"443\u00d73=\"],\n  [\"883\u00d76=\", \"737\u00d76=\"],\n  [\"476\u00d79=\", \"998\u00d78=\"],\n  [\"716\u00d74=\", \"164\u00d72=\"],\n  [\"794\u00d75=\", \"795\u00d77=\"],\n  [\"428\u00d77=\", \"536\u00d74=\"],\n  [\"843\u00d79=\", \"673\u00d72=\"],\n  [\"649\u00d78=\", \"246\u00d79=\"],\n  [\"120\u00d76=\", \"723\u00d79=\"],\n  [\"514\u00d78=\", \"112\u00d77=\"],\n  [\"564\u00d72=\", \"251\u00d75=\"],\n  [\"424\u00d78=\", \"149\u00d72=\"],\n  [\"951\u00d75=\", \"520\u00d79=\"],\n  [\"423\u00d72=\", \"746\u00d77=\"],\n  [\"144\u00d77=\", \"906\u00d77=\"],\n  [\"692\u00d72=\", \"118\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with the\n# updated operands/equation using Find/Replace (wdReplaceOne), which\n# swaps only the matched text and leaves the run's formatting intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"933\u00d73=\", \"361\u00d75=\"),\n  @(\"234\u00d79=\", \"358\u00d79=\"),\n  @(\"137\u00d73=\", \"382\u00d76=\"),\n  @(\"335\u00d77=\", \"732\u00d74=\"),\n  @(\"479\u00d75=\", \"577\u00d78=\"),\n  @(\"550\u00d79=\", \"685\u00d73=\"),\n  @(\"375\u00d76=\", \"219\u00d77=\"),\n  @(\"680\u00d75=\", \"558\u00d79=\"),\n  @(\"275\u00d74=\", \"195\u00d78=\"),\n  @(\"694\u00d77=\", \"443\u00d73=\"),\n  @(\"883\u00d76=\", \"737\u00d76=\"),\n  @(\"476\u00d79=\", \"998\u00d78=\"),\n  @(\"716\u00d74=\", \"164\u00d72=\"),\n  @(\"794\u00d75=\", \"795\u00d77=\"),\n  @(\"428\u00d77=\", \"536\u00d74=\"),\n  @(\"843\u00d79=\", \"673\u00d72=\"),\n  @(\"649\u00d78=\", \"246\u00d79=\"),\n  @(\"120\u00d76=\", \"723\u00d79=\"),\n  @(\"514\u00d78=\", \"112\u00d77=\"),\n  @(\"564\u00d72=\", \"251\u00d75=\"),\n  @(\"424\u00d78=\", \"149\u00d72=\"),\n  @(\"951\u00d75=\", \"520\u00d79=\"),\n  @(\"423\u00d72=\", \"746\u00d77=\"),\n  @(\"144\u00d77=\", \"906\u00d77=\"),\n  @(\"692\u00d72=\", \"118\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
